# [Kadastro App] Yeni kayit eklendi: 2932
# Adds the new record row (Kayit No 2932) to both the master "Kayitlar"
# log sheet and the "Erdemli" unit sheet (a filtered view of the same
# log), appending it right after the existing last row on each sheet.

$wb = $excel.ActiveWorkbook

# --- "Kayitlar" master sheet: new row 20 (after existing row 19) ---
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")

# Force text storage (matches the rest of the sheet, where every column -
# including numeric-looking ones - is stored as text) by switching the
# cell to the "Text" number format before assigning the value, then
# reverting to the workbook's default style so no new persistent
# formatting is introduced.
$wsKayitlar.Cells.Item(20, 1).NumberFormat = "@"
$wsKayitlar.Cells.Item(20, 1).Value = "2932"
$wsKayitlar.Cells.Item(20, 1).Style = "Normal"

$wsKayitlar.Cells.Item(20, 2).NumberFormat = "@"
$wsKayitlar.Cells.Item(20, 2).Value = "2025-09-08"
$wsKayitlar.Cells.Item(20, 2).Style = "Normal"

$wsKayitlar.Cells.Item(20, 3).Value = "Erdemli"

$wsKayitlar.Cells.Item(20, 4).NumberFormat = "@"
$wsKayitlar.Cells.Item(20, 4).Value = "1"
$wsKayitlar.Cells.Item(20, 4).Style = "Normal"

$wsKayitlar.Cells.Item(20, 5).Value = "ÇAP"
$wsKayitlar.Cells.Item(20, 6).Value = "CEMAL TİMUROĞLU (K.Teknisyeni)"

# --- "Erdemli" unit sheet: new row 19 (after existing row 18) ---
$wsErdemli = $wb.Worksheets.Item("Erdemli")

$wsErdemli.Cells.Item(19, 1).NumberFormat = "@"
$wsErdemli.Cells.Item(19, 1).Value = "2932"
$wsErdemli.Cells.Item(19, 1).Style = "Normal"

$wsErdemli.Cells.Item(19, 2).NumberFormat = "@"
$wsErdemli.Cells.Item(19, 2).Value = "2025-09-08"
$wsErdemli.Cells.Item(19, 2).Style = "Normal"

$wsErdemli.Cells.Item(19, 3).Value = "Erdemli"

$wsErdemli.Cells.Item(19, 4).NumberFormat = "@"
$wsErdemli.Cells.Item(19, 4).Value = "1"
$wsErdemli.Cells.Item(19, 4).Style = "Normal"

$wsErdemli.Cells.Item(19, 5).Value = "ÇAP"
$wsErdemli.Cells.Item(19, 6).Value = "CEMAL TİMUROĞLU (K.Teknisyeni)"
